$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header renames: spaces/parentheses/degree marks -> underscores, units stripped ---
$ws.Range("A1").Value  = "Mombo_ShotID"
$ws.Range("G1").Value  = "Ball_mph"
$ws.Range("H1").Value  = "Club_mph"
$ws.Range("I1").Value  = "Smash_Factor"
$ws.Range("J1").Value  = "Carry_yds"
$ws.Range("K1").Value  = "Total_yds"
$ws.Range("L1").Value  = "Roll_yds"
$ws.Range("M1").Value  = "Swing_H"
$ws.Range("N1").Value  = "Spin_rpm"
$ws.Range("O1").Value  = "Height_ft"
$ws.Range("P1").Value  = "Time_s"
$ws.Range("Q1").Value  = "AOA"
$ws.Range("R1").Value  = "Spin_Loft"
$ws.Range("S1").Value  = "Swing_V"
$ws.Range("T1").Value  = "Spin_Axis"
$ws.Range("U1").Value  = "Lateral_yds"
$ws.Range("V1").Value  = "Shot_Type"
$ws.Range("W1").Value  = "FTP"
$ws.Range("X1").Value  = "FTT"
$ws.Range("Y1").Value  = "Dynamic_Loft"
$ws.Range("Z1").Value  = "Club_Path"
$ws.Range("AA1").Value = "Launch_H"
$ws.Range("AB1").Value = "Launch_V"
$ws.Range("AC1").Value = "Low_Point_ftin"
$ws.Range("AD1").Value = "DescentV"
$ws.Range("AE1").Value = "Curve_Dist_yds"
$ws.Range("AF1").Value = "Lateral_Impact_in"
$ws.Range("AG1").Value = "Vertical_Impact_in"
$ws.Range("AJ1").Value = "Unnamed_35"
$ws.Range("AK1").Value = "Unnamed_36"
$ws.Range("AL1").Value = "Unnamed_37"
$ws.Range("AM1").Value = "Unnamed_38"
$ws.Range("AN1").Value = "Unnamed_39"
$ws.Range("AO1").Value = "Unnamed_40"

# --- Data rows 2-6: Spin Axis (T) and Launch H (AA) gain an " L" suffix as text,
#     and Low Point (AC) loses its trailing inch mark and becomes numeric ---
$spinAxis = @{2 = "13.8 L"; 3 = "9.5 L"; 4 = "7.3 L"; 5 = "12.2 L"; 6 = "3.3 L"}
$launchH  = @{2 = "4.9 L";  3 = "4.2 L"; 4 = "4.8 L"; 5 = "6.6 L";  6 = "4.9 L"}
$lowPoint = @{2 = 4.6;      3 = 5.9;     4 = 5.3;     5 = 4.6;     6 = 5.2}

foreach ($r in 2..6) {
    $ws.Range("T$r").Value  = $spinAxis[$r]
    $ws.Range("AA$r").Value = $launchH[$r]
    $ws.Range("AC$r").Value = $lowPoint[$r]
}
